# "Changes to DMD case"
#
# The clinical-exam sheet drops three groups of 4 columns that are no longer
# collected for this case: Hip_add_* (originally A:D), Knee_ext_* (originally
# I:L) and RF_* (originally Q:T). Everything else (Hip_ext_*, Ham_*,
# Gastroc_*, Soleus_*) stays, sliding left to fill the gaps.
#
# Deleting right-to-left means each EntireColumn.Delete() call still targets
# the original column letters, since earlier deletes only affect columns to
# their right (already handled) not columns still to their left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q1:T2").EntireColumn.Delete()   # RF_r_length, RF_l_length, RF_r_stiff, RF_l_stiff
$ws.Range("I1:L2").EntireColumn.Delete()   # Knee_ext_r_ROM, Knee_ext_l_ROM, Knee_ext_r_stiff, Knee_ext_l_stiff
$ws.Range("A1:D2").EntireColumn.Delete()   # Hip_add_r_ROM, Hip_add_l_ROM, Hip_add_r_stiff, Hip_add_l_stiff

# The saved view now has the cursor on the (shifted) Soleus_r_stiff cell.
$ws.Range("O10").Select()

$wb.Save()
